$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of famicon adventure gamebooks to append beneath the existing data.
# Columns: year, japanese, english, publisher, image, product_type
$newRows = @(
    @(1988, "超時空パイレーツ おみそれ3人組の冒険", "Super Space Time Pirates", "Futabasha", "space_time_pirates.jpeg", "gamebook"),
    @(1988, "トキメキハイスクール 恋の学園祭大作戦", "Tokimeki High School", "Futabasha", "tokimeki_high_school.jpeg", "gamebook"),
    @(1988, "霊幻道士 キョンシー大戦争", "Mr. Vampire Jiangshi Great War", "Futabasha", "mr_vampire.jpeg", "gamebook"),
    @(1988, "ガイアの紋章 エルスリード英雄列伝", "Gaia's Coat of Arms", "Futabasha", "gaias_coat_of_arms.jpg", "gamebook"),
    @(1988, "ヤマト魔神伝 サギリ見参! ", "Visit Yamato Genie", "Futabasha", "visit_yamato_genie.jpg", "gamebook"),
    @(1988, "ディープダンジョンIII", "Deep Dungeon III", "Futabasha", "deep_dungeon_3.jpeg", "gamebook"),
    @(1988, "ドラゴンロック 浮遊要塞の死闘 ", "Dragon Rock", "Futabasha", "dragon_rock.jpeg", "gamebook"),
    @(1988, "源平討魔伝 神異妖魔界の変", "Genpei Tomaden", "Futabasha", "genpei_tomaden.jpeg", "gamebook"),
    @(1988, "暗黒要塞ガルディアン オセロ神話の謎", "Dark Fortress Guardian", "Futabasha", "dark_fortress_guardian.jpg", "gamebook"),
    @(1988, "貝獣物語 シェルドラド伝説", "Kaiju Monogatari", "Futabasha", "kaiju_monogatari.jpeg", "gamebook"),
    @(1988, "スペース・ハリアー ホワイトドラゴンの勇者", "Space Harrier", "Futabasha", "space_harrier.jpeg", "gamebook")
)

$startRow = 86
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

$endRow = $startRow + $newRows.Count - 1

# Match the formatting (right border style) already used on column F.
$null = $ws.Range("F85").Copy()
$null = $ws.Range("F$startRow`:F$endRow").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the view to match the post-edit selection/scroll position.
$activeWindow = $excel.ActiveWindow
$null = $activeWindow.SetTopLeftVisibleCell(75, 1)
$activeWindow.ScrollRow = 75
$activeWindow.ScrollColumn = 1
$null = $ws.Range("F85:F$endRow").Select()
